$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.543.78'
$ws.Range("E2").Value = '  +4.00%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.990.24'
$ws.Range("E3").Value = '  +6.21%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9960'
$ws.Range("E4").Value = '  -0.53%  '

# Row 5: XRP
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8177'
$ws.Range("E5").Value = '  +74.46%  '

# Row 6: BNB
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '252.81'
$ws.Range("E6").Value = '  +3.71%  '

# Row 7: USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9970'
$ws.Range("E7").Value = '  -0.44%  '

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3472'
$ws.Range("E8").Value = '  +20.97%  '

# Row 9: Solana
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.06'
$ws.Range("E9").Value = '  +18.90%  '

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06979'
$ws.Range("E10").Value = '  +8.60%  '

# Row 11: Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8410'
$ws.Range("E11").Value = '  +16.75%  '

# Row 12: TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08120'
$ws.Range("E12").Value = '  +4.37%  '

# Row 13: Litecoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '102.43'
$ws.Range("E13").Value = '  +7.67%  '

# Row 14: WrappedEther
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.981.83'
$ws.Range("E14").Value = '  +5.66%  '

# Row 15: Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.509'
$ws.Range("E15").Value = '  +7.50%  '

# Row 16: BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '276.01'
$ws.Range("E16").Value = '  -0.92%  '

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.507.77'
$ws.Range("E17").Value = '  +3.89%  '

# Row 18: Avalanche
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.03'
$ws.Range("E18").Value = '  +8.15%  '

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007903'
$ws.Range("E19").Value = '  +6.87%  '

# Row 20: Uniswap
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.708'
$ws.Range("E20").Value = '  +9.51%  '

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.244.80'
$ws.Range("E21").Value = '  +5.41%  '

# Row 22: Dai
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9981'
$ws.Range("E22").Value = '  -0.33%  '

# Row 23: BinanceUSD
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9946'
$ws.Range("E23").Value = '  -0.65%  '

# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.915'
$ws.Range("E24").Value = '  +10.91%  '

# Row 25: Stellar
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1602'
$ws.Range("E25").Value = '  +67.52%  '

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.700'
$ws.Range("E26").Value = '  +7.62%  '

# Row 27: Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '165.79'
$ws.Range("E27").Value = '  +1.31%  '

# Row 28: EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.81'
$ws.Range("E28").Value = '  +6.33%  '

# Row 29: LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.232'
$ws.Range("E29").Value = '  +18.85%  '

# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.565'
$ws.Range("E30").Value = '  +6.92%  '

# Row 31: Toncoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.356'
$ws.Range("E31").Value = '  +1.90%  '

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.577'
$ws.Range("E32").Value = '  +8.95%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.332'
$ws.Range("E33").Value = '  +6.09%  '

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05226'
$ws.Range("E34").Value = '  +8.92%  '

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.223'
$ws.Range("E35").Value = '  +9.40%  '

# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7494'
$ws.Range("E36").Value = '  +9.51%  '

# Row 37: HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.780'
$ws.Range("E37").Value = '  +2.61%  '

# Row 38: VeChain
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01992'
$ws.Range("E38").Value = '  +6.78%  '

# Row 39: MXToken
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.911'
$ws.Range("E39").Value = '  +3.58%  '

# Row 40: FraxShare
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.611'
$ws.Range("E40").Value = '  +6.40%  '

# Row 41: Aave
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.53'
$ws.Range("E41").Value = '  +5.68%  '

# Row 42: TheSandbox
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4668'
$ws.Range("E42").Value = '  +10.70%  '

# Row 43: RenderToken
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.075'
$ws.Range("E43").Value = '  +7.63%  '

# Row 44: Quant
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '106.17'
$ws.Range("E44").Value = '  +5.48%  '

# Row 45: TrustWalletToken
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8567'
$ws.Range("E45").Value = '  +4.25%  '

# Row 46: PaxDollar
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9978'
$ws.Range("E46").Value = '  -0.25%  '

# Row 47: EnergySwap
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.934'
$ws.Range("E47").Value = '  +3.61%  '

# Row 48: Aptos
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.513'
$ws.Range("E48").Value = '  +8.76%  '

# Row 49: Elrond
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.62'
$ws.Range("E49").Value = '  +4.61%  '

# Row 50: Decentraland
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4274'
$ws.Range("E50").Value = '  +9.90%  '

# Row 51: Maker
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '933.01'
$ws.Range("E51").Value = '  +4.00%  '
